$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "accuracy" column (E), replacing the old single G5/H5 summary cell ---
$ws.Range("E1").Value = "accuracy"

# E2 is a standalone formula; E3:E18 form one shared-formula group (mirrors column D's pattern)
$ws.Range("E2").Formula = "=1 - (D2/A2)"
$ws.Range("E3:E18").Formula = "=1 - (D3/A3)"

# Give E2:E18 the "Percent" based style used for the accuracy figures
$ws.Range("E2:E18").Style = "Percent"

# --- Remove the old single accuracy cell (G5 label + H5 formula) ---
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()

# --- Summary rows: AVERAGE / STDEV / MAX / MIN across A:E ---
$ws.Range("A21:B21").Formula = "=AVERAGE(A2:A18)"
$ws.Range("C21").Formula = "=AVERAGE(C2:C18)"
$ws.Range("D21:E21").Formula = "=AVERAGE(D2:D18)"

$ws.Range("A22:B22").Formula = "=STDEV(A2:A18)"
$ws.Range("C22").Formula = "=STDEV(C2:C18)"
$ws.Range("D22:E22").Formula = "=STDEV(D2:D18)"

$ws.Range("A23").Formula = "=MAX(A2:A18)"
$ws.Range("B23:E23").Formula = "=MAX(B2:B18)"

$ws.Range("A24").Formula = "=MIN(A2:A18)"
$ws.Range("B24:E24").Formula = "=MIN(B2:B18)"

# Rows 21/22 (A:D) carry the same integer style as the existing B21 seed cell;
# E21/E22 carry the percent style like the rest of column E
$ws.Range("A21:D22").Style = "Comma [0]"
$ws.Range("A21:D22").NumberFormat = "General"
$ws.Range("E21:E24").Style = "Percent"

# --- Column H is no longer a data column; widen it for the (now blank) cell ---
$ws.Columns.Item(8).ColumnWidth = 22.43

# --- Selection moves to the new STDEV/accuracy cell ---
$ws.Range("E22").Select()
